# KP-11725 D: Extension of questionnaire's translation files
# Insert a new "Variable" column right after the "Entity Id" column so the
# translations sheet also carries each question/entity's variable name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "Variable"
$ws.Range("B2").Value = "e1"

$ws.Columns("B:B").ColumnWidth = 18.67

$ws.Range("B3").Select() | Out-Null
